# Refresh the crypto price/volume snapshot (cells B/C/D/E, rows 2-51)
# to match the latest scrape, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '63.434.65'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  +1.09%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '3.410.00'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  +1.92%  '; ForceText = $false }
    @{ Cell = 'E4'; Value = '  -0.04%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '568.27'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  +1.06%  '; ForceText = $false }
    @{ Cell = 'D6'; Value = '155.78'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  +2.17%  '; ForceText = $false }
    @{ Cell = 'E7'; Value = '  +0.01%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '3.409.33'; ForceText = $false }
    @{ Cell = 'E8'; Value = '  +1.71%  '; ForceText = $false }
    @{ Cell = 'E9'; Value = '  +2.44%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '7.41'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  -0.52%  '; ForceText = $false }
    @{ Cell = 'E11'; Value = '  +3.63%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '0.433'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  -0.49%  '; ForceText = $false }
    @{ Cell = 'D13'; Value = '3.985.58'; ForceText = $false }
    @{ Cell = 'E13'; Value = '  +1.62%  '; ForceText = $false }
    @{ Cell = 'E14'; Value = '  -3.00%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '0.0000193'; ForceText = $true }
    @{ Cell = 'E15'; Value = '  +7.74%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '27.22'; ForceText = $true }
    @{ Cell = 'E16'; Value = '  +1.15%  '; ForceText = $false }
    @{ Cell = 'D17'; Value = '63.479.34'; ForceText = $false }
    @{ Cell = 'E17'; Value = '  +1.18%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '3.391.27'; ForceText = $false }
    @{ Cell = 'E18'; Value = '  +2.22%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '6.26'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  -1.36%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '14.08'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  +1.82%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '380.34'; ForceText = $true }
    @{ Cell = 'E21'; Value = '  -1.03%  '; ForceText = $false }
    @{ Cell = 'D22'; Value = '8.07'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  -3.47%  '; ForceText = $false }
    @{ Cell = 'E23'; Value = '  +0.59%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '71.61'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  +2.07%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '0.529'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  -1.19%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '0.0000121'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  +27.28%  '; ForceText = $false }
    @{ Cell = 'E27'; Value = '  +5.94%  '; ForceText = $false }
    @{ Cell = 'E28'; Value = '  +0.08%  '; ForceText = $false }
    @{ Cell = 'D30'; Value = '6.05'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  +8.55%  '; ForceText = $false }
    @{ Cell = 'D31'; Value = '1.37'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  +4.78%  '; ForceText = $false }
    @{ Cell = 'E32'; Value = '  +1.04%  '; ForceText = $false }
    @{ Cell = 'B33'; Value = 'EthereumClassic'; ForceText = $false }
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; ForceText = $false }
    @{ Cell = 'D33'; Value = '23.24'; ForceText = $true }
    @{ Cell = 'E33'; Value = '  +1.44%  '; ForceText = $false }
    @{ Cell = 'B34'; Value = 'RenderToken'; ForceText = $false }
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; ForceText = $false }
    @{ Cell = 'D34'; Value = '6.41'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  -2.57%  '; ForceText = $false }
    @{ Cell = 'E35'; Value = '  +0.00%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '6.79'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  +1.46%  '; ForceText = $false }
    @{ Cell = 'D37'; Value = '159.55'; ForceText = $true }
    @{ Cell = 'E37'; Value = '  -0.28%  '; ForceText = $false }
    @{ Cell = 'D38'; Value = '1.45'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  -1.44%  '; ForceText = $false }
    @{ Cell = 'D39'; Value = '2.967.54'; ForceText = $false }
    @{ Cell = 'E39'; Value = '  +5.62%  '; ForceText = $false }
    @{ Cell = 'B40'; Value = 'Hedera'; ForceText = $false }
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; ForceText = $false }
    @{ Cell = 'D40'; Value = '0.0757'; ForceText = $true }
    @{ Cell = 'E40'; Value = '  +2.48%  '; ForceText = $false }
    @{ Cell = 'B41'; Value = 'EnergySwap'; ForceText = $false }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; ForceText = $false }
    @{ Cell = 'D41'; Value = '26.99'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  +0.05%  '; ForceText = $false }
    @{ Cell = 'E42'; Value = '  -3.54%  '; ForceText = $false }
    @{ Cell = 'D43'; Value = '0.0316'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  +0.90%  '; ForceText = $false }
    @{ Cell = 'D44'; Value = '41.78'; ForceText = $true }
    @{ Cell = 'E44'; Value = '  +3.01%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '0.765'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  +2.67%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '4.33'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  +1.73%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '23.31'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  +6.24%  '; ForceText = $false }
    @{ Cell = 'D48'; Value = '1.08'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  +3.64%  '; ForceText = $false }
    @{ Cell = 'D49'; Value = '2.22'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  +23.60%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '0.835'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  +4.62%  '; ForceText = $false }
    @{ Cell = 'D51'; Value = '6.35'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  +0.82%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # These D-column prices are plain text in the sheet (e.g. "7.41"), but
        # Excel auto-converts numeric-looking text to a Number on assignment.
        # Force the Text format first, write the value, then drop back to the
        # workbook default style so only the cell content changes.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
